$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in the header cell
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 15:35"

# Full target data for the provincias table (rows 4-57): Ciudad, Casos totales, Casos activos, Recuperados, Muertes
$data = @(
    @(4, 'Madrid', 66338, 40698, 16793, 8847),
    @(5, 'Cataluña', 55824, 26172, 23708, 5944),
    @(6, 'Castilla y Leon', 18470, 8716, 7808, 1946),
    @(7, 'Castilla-La Mancha', 16618, 6378, 7347, 2893),
    @(8, 'Pais Vasco', 13156, 14646, 0, 1418),
    @(9, 'Andalucia', 12450, 10611, 481, 1358),
    @(10, 'Bizkaia/Vizcaya', 10332, 7124, 4423, 815),
    @(11, 'Galicia', 9041, 8409, 28, 604),
    @(12, 'Ciudad Real', 6464, 2030, 3368, 1066),
    @(13, 'Valencia/Valencia', 5609, 4907, 2767, 693),
    @(14, 'Aragon', 5478, 3727, 913, 838),
    @(15, 'Zaragoza', 5287, 2820, 1812, 655),
    @(16, 'Navarra', 5148, 3751, 894, 503),
    @(17, 'Araba/Alava', 4868, 7124, 4423, 356),
    @(18, 'Valladolid', 4393, 1577, 2454, 362),
    @(19, 'Salamanca', 4152, 1181, 2613, 358),
    @(20, 'La Rioja', 4024, 3048, 627, 349),
    @(21, 'Toledo', 3872, 1992, 1124, 756),
    @(22, 'Alacant/Alicante', 3794, 3637, 1938, 484),
    @(23, 'Albacete', 3775, 1389, 1872, 514),
    @(24, 'Leon', 3569, 1629, 1533, 407),
    @(25, 'Segovia', 3413, 868, 2344, 201),
    @(26, 'Gipuzkoa/Guipuzcoa', 3116, 7124, 4423, 283),
    @(27, 'Extremadura', 2919, 2422, 10, 487),
    @(28, 'Malaga', 2758, 2295, 185, 278),
    @(29, 'Burgos', 2746, 901, 1640, 205),
    @(30, 'Sevilla', 2423, 1840, 308, 275),
    @(31, 'Granada', 2413, 2520, 0, 279),
    @(32, 'Asturias', 2366, 1061, 990, 315),
    @(33, 'Soria', 2290, 397, 1774, 119),
    @(34, 'Gran Canaria', 2289, 1524, 614, 151),
    @(35, 'Tenerife', 2280, 1506, 623, 151),
    @(36, 'Cantabria', 2246, 1981, 62, 203),
    @(37, 'Caceres', 1973, 1505, 66, 402),
    @(38, 'A Coruña', 1969, 333, 1788, 67),
    @(39, 'Avila', 1935, 623, 1179, 133),
    @(40, 'Pontevedra', 1536, 333, 1411, 30),
    @(41, 'Murcia', 1508, 1782, 0, 139),
    @(42, 'Castello/Castellon', 1486, 1363, 699, 207),
    @(43, 'Jaen', 1387, 1171, 41, 175),
    @(44, 'Cordoba', 1331, 1350, 0, 106),
    @(45, 'Guadalajara', 1266, 371, 644, 251),
    @(46, 'Cuenca', 1241, 596, 339, 306),
    @(47, 'Cadiz', 1240, 560, 535, 145),
    @(48, 'Palencia', 1205, 333, 789, 83),
    @(49, 'Huesca', 1115, 472, 544, 99),
    @(50, 'Zamora', 993, 322, 586, 85),
    @(51, 'Badajoz', 962, 1082, 0, 92),
    @(52, 'Ourense', 751, 333, 660, 22),
    @(53, 'Teruel', 664, 378, 203, 83),
    @(54, 'Lugo', 586, 333, 520, 11),
    @(55, 'Almeria', 498, 484, 0, 52),
    @(56, 'Huelva', 400, 391, 0, 48),
    @(57, 'Mallorca', 210, 18, 194, 12)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
